# Daily attendance processing - 2026-01-15 19:40:13
# Reverses the order of the comma-separated "Recorded By" entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 157) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value()

    if ($null -eq $text) { continue }
    $text = [string]$text

    if ($text -match ",") {
        $parts = $text -split ","
        $n = $parts.Count
        $rev = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $rev += $parts[$i].Trim()
        }
        $cell.Value = ($rev -join ", ")
    }
}
